$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the wording of the task description in A7 (DB connectie "doet moeilijk" -> "bestaat nog niet")
$ws.Range("A7").Value = "Implementatie van Repositories, services en controllers zodat personen kunnen inloggen op de correcte manier (nog te testen, want connectie met DB bestaat nog niet)"

# 2) Change date-cell alignment (B8:B11) from general to right-aligned
$ws.Range("B8").HorizontalAlignment = -4152
$ws.Range("B9").HorizontalAlignment = -4152
$ws.Range("B10").HorizontalAlignment = -4152
$ws.Range("B11").HorizontalAlignment = -4152

# 3) Fill in the previously-empty row 12 with a new log entry
$ws.Range("A12").Value = "Creeeren van index pagina voor de klant. Deze kan nu al de beschreven acties doen."
$ws.Range("B12").Value = "21/12/2024"
$ws.Range("B12").HorizontalAlignment = -4152
$ws.Range("C12").Value = 3.5
$ws.Range("D12").Value = 0
$ws.Rows.Item(12).RowHeight = 58.2

# 4) Move the active selection to D13
$ws.Range("D13").Select()

Write-Output "done"
